# SprintHours.xlsx edit script
# Applies the Sprint 5 ("Aaron Mulligan" sprint-5 entries) hour updates,
# clears a couple of stray numeric/empty cells, adds the two new
# "1 hour " / "3 hours" notes, deletes the old blank spacer rows 45-48
# and records the "Change the anniversary picture" hours (1 hour) in
# the new blog-update/fuse-box Sprint table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: add a styled (but empty) cell at E12, matching its neighbours ---
$ws.Range("E12").Style = $ws.Range("D12").Style

# --- Row 25: drop the stray G25 total (5) ---
$ws.Range("G25").Clear()

# --- Row 28 (Aaron Mulligan / puzzle memories / 2 hours): record Sprint total ---
$ws.Range("F28").Style = $ws.Range("E28").Style
$ws.Range("G28").Style = $ws.Range("E28").Style
$ws.Range("G28").Value = 6

# --- Row 29: doctor's letter task now logged at 1 hour ---
# (write this BEFORE row 44 so the new shared strings are appended in
#  "1 hour " (36) then "3 hours" (37) order, matching the workbook)
$ws.Range("E29").Value = "1 hour "

# --- Row 36: drop the stray empty G36 cell ---
$ws.Range("G36").Clear()

# --- Row 41 (Aaron Mulligan / presentation & practice): 2 hours, Sprint total 8 ---
$ws.Range("E41").Value = "2 hours"
$ws.Range("F41").Style = $ws.Range("E41").Style
$ws.Range("G41").Style = $ws.Range("E41").Style
$ws.Range("G41").Value = 8

# --- Row 42 (Finalize the objects to go in specific rooms): 1 hour ---
$ws.Range("E42").Value = "1 hour"

# --- Row 43 (design document for puzzle mechanic): 2 hours ---
$ws.Range("E43").Value = "2 hours"

# --- Row 44 (fuse box management mechanic design): 3 hours (new string) ---
$ws.Range("E44").Value = "3 hours"

# --- Remove the 4 blank spacer rows (45:48) above the blog-update table ---
$ws.Rows("45:48").Delete()

# --- New row 47 (old 51): Sprint number bumped from 4 to 5 ---
$ws.Range("F47").Value = 5

# --- New row 54 (old 58, "Change the anniversary picture"): log 1 hour ---
$ws.Range("E54").Value = 1
